$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reverse the "Periodo Mora" (E) and "Valor Mora" (F) columns for rows 16-22
# (previous account statements removed, new ones added - data re-ordered)
$ws.Range("E16").Value = "2108"
$ws.Range("F16").Value = 35129

$ws.Range("E17").Value = "2107"
$ws.Range("F17").Value = 36341

$ws.Range("E18").Value = "2106"
$ws.Range("F18").Value = 36341

$ws.Range("E19").Value = "2105"
$ws.Range("F19").Value = 36341

$ws.Range("E20").Value = "2104"
$ws.Range("F20").Value = 36341

$ws.Range("E21").Value = "2103"
$ws.Range("F21").Value = 36341

$ws.Range("E22").Value = "2102"
$ws.Range("F22").Value = 36341
